$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.172059178352356
$ws.Range("B1").Value = 2.122202634811401
$ws.Range("C1").Value = 3.150188684463501
$ws.Range("D1").Value = 0.435005247592926
$ws.Range("E1").Value = 1.373210787773132
